$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AYKO")

# Insert a new row at position 43, shifting existing rows 43.. down by one.
$ws.Rows("43").Insert()

# Populate the newly inserted row 43 with the new record's data.
# Leading "'" forces text storage for the numeric-looking fields so the
# cell type matches the rest of the (text) data columns (A-L), mirroring
# the workbook's existing convention of storing these as inline strings.
$ws.Range("A43").Value = "'5521"
$ws.Range("B43").Value = "'4/8/2025"
$ws.Range("C43").Value = "EL PEREGRINO 3115"
$ws.Range("D43").Value = "'11"
$ws.Range("E43").Value = "'804569000"
$ws.Range("F43").Value = "AYKO"
$ws.Range("G43").Value = "Pendiente"
$ws.Range("H43").Value = "Volvio a ingresar se inclino el poste - caso 6316"
$ws.Range("I43").Value = "'1"
$ws.Range("J43").Value = "Aplomo"
$ws.Range("K43").Value = "Sin equipos"
$ws.Range("L43").Value = "Poste"
$ws.Range("M43").Value = -58.485232
$ws.Range("N43").Value = -34.611573
